# Scheduled market-data refresh for the Lamia_Profits leve-crafting workbook.
# Re-pastes updated Universalis average-price snapshots (columns H:N --
# currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) for the rows whose
# backing market data moved since the last run. Values only -- no formulas
# live in this range, so each cell is written directly with the refreshed
# number pulled from the upstream price snapshot.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 7583.2607
$ws.Range("I62").Value = 6339.615
$ws.Range("K62").Value = 6339.615
$ws.Range("M62").Value = -5715.615
# Row 65
$ws.Range("H65").Value = 7583.2607
$ws.Range("I65").Value = 6339.615
$ws.Range("K65").Value = 31698.075
$ws.Range("M65").Value = -28578.075
# Row 70
$ws.Range("H70").Value = 11116914
$ws.Range("I70").Value = 5068.125
$ws.Range("K70").Value = 15204.375
$ws.Range("M70").Value = -14934.375
# Row 73
$ws.Range("H73").Value = 11116914
$ws.Range("I73").Value = 5068.125
$ws.Range("K73").Value = 15204.375
$ws.Range("M73").Value = -14268.375
# Row 74
$ws.Range("H74").Value = 10257.667
$ws.Range("I74").Value = 9498.799999999999
$ws.Range("K74").Value = 9498.799999999999
$ws.Range("M74").Value = -8562.799999999999
# Row 77
$ws.Range("H77").Value = 10257.667
$ws.Range("I77").Value = 9498.799999999999
$ws.Range("K77").Value = 47494
$ws.Range("M77").Value = -42814
# Row 118
$ws.Range("H118").Value = 516.7143
$ws.Range("I118").Value = 516.7143
$ws.Range("K118").Value = 1550.1429
$ws.Range("M118").Value = 106.8571000000002
# Row 131
$ws.Range("H131").Value = 143415.17
$ws.Range("I131").Value = 2406.7273
$ws.Range("K131").Value = 7220.1819
$ws.Range("M131").Value = -2180.1819
# Row 132
$ws.Range("H132").Value = 1253.0588
$ws.Range("I132").Value = 1080.4286
$ws.Range("K132").Value = 3241.2858
$ws.Range("M132").Value = -711.2857999999997

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5383.4136
$ws.Range("I32").Value = 4683.8
$ws.Range("K32").Value = 4683.8
$ws.Range("M32").Value = -4396.8
# Row 74
$ws.Range("H74").Value = 55561510
$ws.Range("I74").Value = 66672396
$ws.Range("J74").Value = 7100
$ws.Range("K74").Value = 66672396
$ws.Range("L74").Value = 7100
$ws.Range("M74").Value = -66671522
$ws.Range("N74").Value = -8848
# Row 77
$ws.Range("H77").Value = 55561510
$ws.Range("I77").Value = 66672396
$ws.Range("J77").Value = 7100
$ws.Range("K77").Value = 333361980
$ws.Range("L77").Value = 35500
$ws.Range("M77").Value = -333357612
$ws.Range("N77").Value = -44236
# Row 96
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
# Row 97
$ws.Range("H97").Value = 737
$ws.Range("I97").Value = 913.1
$ws.Range("K97").Value = 913.1
$ws.Range("M97").Value = -417.1
# Row 122
$ws.Range("H122").Value = 3711.25
$ws.Range("I122").Value = 2768.75
$ws.Range("K122").Value = 8306.25
$ws.Range("M122").Value = -5856.25
# Row 132
$ws.Range("H132").Value = 3586.8125
$ws.Range("I132").Value = 2558.3333
$ws.Range("J132").Value = 19014
$ws.Range("K132").Value = 7674.999899999999
$ws.Range("L132").Value = 57042
$ws.Range("M132").Value = -5144.999899999999
$ws.Range("N132").Value = -62102
# Row 134
$ws.Range("H134").Value = 122945
$ws.Range("J134").Value = 122945
$ws.Range("L134").Value = 122945
$ws.Range("N134").Value = -133085

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1234.2
$ws.Range("I94").Value = 987.25
$ws.Range("K94").Value = 987.25
$ws.Range("M94").Value = -536.25
# Row 105
$ws.Range("H105").Value = 12491.167
$ws.Range("I105").Value = 17186.846
$ws.Range("K105").Value = 17186.846
$ws.Range("M105").Value = -15439.846

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 1571.75
$ws.Range("I22").Value = 1703.3572
$ws.Range("J22").Value = 1469.3889
$ws.Range("K22").Value = 1703.3572
$ws.Range("L22").Value = 1469.3889
$ws.Range("M22").Value = -1353.3572
$ws.Range("N22").Value = -2169.3889
# Row 31
$ws.Range("H31").Value = 35999.793
$ws.Range("I31").Value = 2955
$ws.Range("K31").Value = 2955
$ws.Range("M31").Value = -2660
# Row 34
$ws.Range("H34").Value = 35999.793
$ws.Range("I34").Value = 2955
$ws.Range("K34").Value = 2955
$ws.Range("M34").Value = -2753
# Row 58
$ws.Range("H58").Value = 5017.2085
$ws.Range("I58").Value = 3303.125
$ws.Range("J58").Value = 8445.375
$ws.Range("K58").Value = 3303.125
$ws.Range("L58").Value = 8445.375
$ws.Range("M58").Value = -3100.125
$ws.Range("N58").Value = -8851.375
# Row 62
$ws.Range("H62").Value = 5682.4443
$ws.Range("I62").Value = 2222.5
$ws.Range("J62").Value = 12602.333
$ws.Range("K62").Value = 2222.5
$ws.Range("L62").Value = 12602.333
$ws.Range("M62").Value = -1598.5
$ws.Range("N62").Value = -13850.333
# Row 65
$ws.Range("H65").Value = 5682.4443
$ws.Range("I65").Value = 2222.5
$ws.Range("J65").Value = 12602.333
$ws.Range("K65").Value = 11112.5
$ws.Range("L65").Value = 63011.665
$ws.Range("M65").Value = -7992.5
$ws.Range("N65").Value = -69251.66500000001
# Row 107
$ws.Range("H107").Value = 1088.2632
$ws.Range("I107").Value = 872.2
$ws.Range("K107").Value = 872.2
$ws.Range("M107").Value = 1047.8
# Row 122
$ws.Range("H122").Value = 8828.454
$ws.Range("I122").Value = 4266.7144
$ws.Range("K122").Value = 12800.1432
$ws.Range("M122").Value = -10350.1432
# Row 125
$ws.Range("H125").Value = 89546.664
$ws.Range("J125").Value = 89546.664
$ws.Range("L125").Value = 89546.664
$ws.Range("N125").Value = -94466.664
# Row 132
$ws.Range("H132").Value = 6386.4614
$ws.Range("I132").Value = 5001.25
$ws.Range("K132").Value = 15003.75
$ws.Range("M132").Value = -12473.75
# Row 136
$ws.Range("H136").Value = 5017.2085
$ws.Range("I136").Value = 3303.125
$ws.Range("J136").Value = 8445.375
$ws.Range("K136").Value = 9909.375
$ws.Range("L136").Value = 25336.125
$ws.Range("M136").Value = -7359.375
$ws.Range("N136").Value = -30436.125

$ws = $wb.Worksheets.Item("CUL")
# Row 55
$ws.Range("H55").Value = 300
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
# Row 109
$ws.Range("H109").Value = 1752.6875
$ws.Range("I109").Value = 1715.8889
$ws.Range("K109").Value = 5147.6667
$ws.Range("M109").Value = -4107.6667
# Row 139
$ws.Range("H139").Value = 4311.1177
$ws.Range("I139").Value = 1362.909
$ws.Range("J139").Value = 9716.166999999999
$ws.Range("K139").Value = 4088.727
$ws.Range("L139").Value = 29148.501
$ws.Range("M139").Value = 1051.273
$ws.Range("N139").Value = -39428.501

$ws = $wb.Worksheets.Item("GSM")
# Row 7
$ws.Range("H7").Value = 18000
$ws.Range("J7").Value = 18000
$ws.Range("L7").Value = 18000
$ws.Range("N7").Value = -18224
# Row 8
$ws.Range("H8").Value = 18000
$ws.Range("J8").Value = 18000
$ws.Range("L8").Value = 18000
$ws.Range("N8").Value = -18278
# Row 123
$ws.Range("H123").Value = 33326
$ws.Range("J123").Value = 33326
$ws.Range("L123").Value = 33326
$ws.Range("N123").Value = -38226
# Row 132
$ws.Range("H132").Value = 120285.78
$ws.Range("I132").Value = 149220.72
$ws.Range("J132").Value = 19013.5
$ws.Range("K132").Value = 447662.16
$ws.Range("L132").Value = 57040.5
$ws.Range("M132").Value = -445132.16
$ws.Range("N132").Value = -62100.5

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 16204.6
$ws.Range("I40").Value = 16190.167
$ws.Range("K40").Value = 16190.167
$ws.Range("M40").Value = -16054.167
# Row 132
$ws.Range("H132").Value = 5800
$ws.Range("J132").Value = 9402
$ws.Range("L132").Value = 28206
$ws.Range("N132").Value = -33266
# Row 136
$ws.Range("H136").Value = 8404
$ws.Range("I136").Value = 4257.5
$ws.Range("K136").Value = 12772.5
$ws.Range("M136").Value = -10222.5

$ws = $wb.Worksheets.Item("WVR")
# Row 3
$ws.Range("H3").Value = 5001
$ws.Range("I3").Value = 5001.5
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 5001.5
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = -4887.5
$ws.Range("N3").Value = -5228
# Row 122
$ws.Range("H122").Value = 5966.3335
$ws.Range("I122").Value = 4336.5
$ws.Range("K122").Value = 13009.5
$ws.Range("M122").Value = -10559.5
# Row 128
$ws.Range("H128").Value = 50000
$ws.Range("J128").Value = 50000
$ws.Range("L128").Value = 50000
$ws.Range("N128").Value = -59960
# Row 132
$ws.Range("H132").Value = 6263.4287
$ws.Range("I132").Value = 6018.5
$ws.Range("J132").Value = 6018.5
$ws.Range("K132").Value = 18055.5
$ws.Range("M132").Value = -15525.5

